$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New professor-title rows to append below the existing data (rows 2-9),
# matching columns: title (A), count (B), year (C, stored as text), century (D)
$newRows = @(
    @("Sir", 1, "1931", 20),
    @("Sir", 1, "1932", 20),
    @("Dame", 1, "1937", 20),
    @("Sir", 1, "1941", 20),
    @("Sir", 1, "1943", 20),
    @("Sir", 1, "1944", 20),
    @("Jhr.", 1, "1946", 20)
)

$startRow = 10
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]

    # Year must be stored as text (like the existing rows), not auto-converted
    # to a number by Excel's smart entry. Temporarily force a text format,
    # assign the value, then reset the style so no extra formatting lingers.
    $yearCell = $ws.Cells.Item($r, 3)
    $yearCell.NumberFormat = "@"
    $yearCell.Value = $row[2]
    $yearCell.Style = "Normal"

    $ws.Cells.Item($r, 4).Value = $row[3]
}
